# feat: add 2022-Q3 data
#
# 1. Duplicate the "2021-Q3" sheet (our best template: same columns/styles)
#    and insert the copy immediately before it, then rename to "2022-Q3"
#    and overwrite its data with the new quarter's figures.
# 2. Insert a new row at the top of the "总计" (total) summary sheet for
#    the new quarter, pushing the older rows down, and renumber the index
#    column so it stays a plain 0..n sequence.

$wb = $excel.ActiveWorkbook

# --- Step 1: create & populate the "2022-Q3" worksheet -------------------

$templateSheet = $wb.Worksheets.Item("2021-Q3")
$templateSheet.Copy($templateSheet)

$newSheet = $wb.Worksheets.Item(2)
$newSheet.Name = "2022-Q3"

# Header row: only the fund-size column label changed.
$newSheet.Range("D1").Value = "基金规模"

# Data row - keep these as literal text (matching the source formatting)
# rather than letting plain numeric-looking strings get auto-converted to
# numbers. A leading apostrophe forces text entry without touching the
# cell's number format (same trick as typing '5.40 directly into Excel).
$newSheet.Range("D2").Value = "'5.40"
$newSheet.Range("E2").Value = "'99.62"
$newSheet.Range("F2").Value = "'2.51"
$newSheet.Range("G2").Value = "'0.1355"
$newSheet.Range("H2").Value = 10

# --- Step 2: add the new row to the "总计" summary sheet ------------------

$totalSheet = $wb.Worksheets.Item("总计")
$totalSheet.Rows.Item(2).Insert()

# The freshly inserted row silently inherits the formatting of the row
# above it (the bold/bordered header); the source file has no explicit
# style on these data cells, so strip it back off.
$totalSheet.Range("B2:D2").ClearFormats()

$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q3"
$totalSheet.Range("C2").Value = 1
$totalSheet.Range("D2").Value = 0.14

# Copy the style used by the other index/date cells onto the new row so it
# matches the rest of the table (border + centered + bold).
$totalSheet.Range("A3").Copy()
$totalSheet.Range("A2").PasteSpecial(-4122)

# Renumber the index column (A) to stay a plain 0..n sequence after the
# insert.
$totalSheet.Range("A3").Value = 1
$totalSheet.Range("A4").Value = 2
$totalSheet.Range("A5").Value = 3

$excel.CutCopyMode = 0
